$wb = $excel.ActiveWorkbook

$wsMeans = $wb.Worksheets.Item("Means")
$wsSD    = $wb.Worksheets.Item("Standard Deviations")

# ---------------------------------------------------------------------------
# New column headers (F = "5 miles", G = "10 miles")
# ---------------------------------------------------------------------------
$wsMeans.Range("F1").Value = "Within 5 miles of HFC production facility"
$wsMeans.Range("G1").Value = "Within 10 miles of HFC production facility"

$wsSD.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$wsSD.Range("G1").Value = "Within 10 mile of HFC production facility SD"

# ---------------------------------------------------------------------------
# Means sheet - new column F (5 miles) and G (10 miles) values
# ---------------------------------------------------------------------------
$meansF = @(72, 22, 6.1, 7.1, 56, 9.2, 7.9, 35, 0.41)
$meansG = @(76, 18, 6.6, 5.9, 61, 8.2, 7.2, 34, 0.4)

for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $wsMeans.Cells.Item($row, 6).Value = $meansF[$i]
    $wsMeans.Cells.Item($row, 7).Value = $meansG[$i]
}

# Updated existing Means values for row 9 (Total Cancer Risk) and row 10 (Total Respiratory)
$wsMeans.Range("B9").Value = 29
$wsMeans.Range("C9").Value = 26
$wsMeans.Range("D9").Value = 45
$wsMeans.Range("E9").Value = 40

$wsMeans.Range("B10").Value = 0.37
$wsMeans.Range("C10").Value = 0.32
$wsMeans.Range("D10").Value = 0.48
$wsMeans.Range("E10").Value = 0.42

# ---------------------------------------------------------------------------
# Standard Deviations sheet - new column F (5 miles) and G (10 miles) values
# ---------------------------------------------------------------------------
$sdF = @(30, 30, 7.8, 13, 25, 11, 11, 19, 0.19)
$sdG = @(28, 27, 8.8, 11, 28, 11, 10, 15, 0.14)

for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $wsSD.Cells.Item($row, 6).Value = $sdF[$i]
    $wsSD.Cells.Item($row, 7).Value = $sdG[$i]
}

# Updated existing SD values for row 9 and row 10
$wsSD.Range("B9").Value = 10
$wsSD.Range("C9").Value = 8.6
$wsSD.Range("E9").Value = 28

$wsSD.Range("B10").Value = 0.14
$wsSD.Range("D10").Value = 0.43
$wsSD.Range("E10").Value = 0.24
